{"js": "// Group 18 Contract \u2014 three textual fixes:\n//  1) Item 2: the \"late comer\" sentence is re-typed as a single run\n//     (no visible text change, just a run clean-up in the source doc).\n//  2) Item 6: \"UCI's rules\" -> \"UNCC's rules\".\n//  3) Item 9: the \"timekeeper's\" sentence is re-typed as a single run\n//     (no visible text change, just a run clean-up in the source doc).\n\nconst body = context.document.body;\n\n// --- 1) Paragraph 2 ------------------------------------------------------\nconst p2Text =\n  \"2. In the event that a group member is less than five minutes late, \" +\n  \"s/he may quietly join the group without disrupting it to ask what \" +\n  \"s/he missed. It is optional for the group members to fill in the \" +\n  \"late comer.\";\nconst p2Results = body.search(p2Text, { matchCase: true });\np2Results.load(\"text\");\nawait context.sync();\nif (p2Results.items.length > 0) {\n  p2Results.items[0].insertText(p2Text, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2) Paragraph 6: UCI -> UNCC ------------------------------------------\nconst uciResults = body.search(\"UCI\", { matchCase: true });\nuciResults.load(\"text\");\nawait context.sync();\nif (uciResults.items.length > 0) {\n  uciResults.items[0].insertText(\"UNCC\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 3) Paragraph 9 -------------------------------------------------------\nconst p9Text =\n  \"9. In the event that a group member or members are dominating the \" +\n  \"group, it\\u2019s the timekeeper\\u2019s job to politely interrupt them \" +\n  \"(this is when you can interrupt) and ask that someone else speak. \";\nconst p9Results = body.search(p9Text, { matchCase: true });\np9Results.load(\"text\");\nawait context.sync();\nif (p9Results.items.length > 0) {\n  p9Results.items[0].insertText(p9Text, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Group 18 Contract \u2014 three textual fixes:\n#  1) Item 2: the \"late comer\" sentence is re-typed as a single run\n#     (no visible text change, just a run clean-up in the source doc).\n#  2) Item 6: \"UCI's rules\" -> \"UNCC's rules\".\n#  3) Item 9: the \"timekeeper's\" sentence is re-typed as a single run\n#     (no visible text change, just a run clean-up in the source doc).\n\n$d = $word.ActiveDocument\n\n# --- 1) Paragraph 2 --------------------------------------------------------\n$p2Text = \"2. In the event that a group member is less than five minutes late, s/he may quietly join the group without disrupting it to ask what s/he missed. It is optional for the group members to fill in the late comer.\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $p2Text\n$find.Replacement.Text = $p2Text\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# --- 2) Paragraph 6: UCI -> UNCC --------------------------------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"UCI\"\n$find2.Replacement.Text = \"UNCC\"\n$find2.Execute($find2.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\n# --- 3) Paragraph 9 ----------------------------------------------------------\n$p9Text = \"9. In the event that a group member or members are dominating the group, it\u2019s the timekeeper\u2019s job to politely interrupt them (this is when you can interrupt) and ask that someone else speak. \"\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = $p9Text\n$find3.Replacement.Text = $p9Text\n$find3.Execute($find3.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2)\n"}
